$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header date labels (October -> November) across all five header blocks ---
# NumberFormat is temporarily forced to text ("@") so Excel does not auto-convert
# the "November 2016" / "November 2015" strings into date serial values; the
# original numeric display format is restored immediately afterwards.
$headerRows = @(7,18,29,40,51)
$map2016 = @("B","D","F","H","J")
$map2015 = @("C","E","G","I","K")
foreach ($r in $headerRows) {
    foreach ($col in $map2016) {
        $addr = "$col$r"
        $origFmt = $ws.Range($addr).NumberFormat
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = "November 2016"
        $ws.Range($addr).NumberFormat = $origFmt
    }
    foreach ($col in $map2015) {
        $addr = "$col$r"
        $origFmt = $ws.Range($addr).NumberFormat
        $ws.Range($addr).NumberFormat = "@"
        $ws.Range($addr).Value = "November 2015"
        $ws.Range($addr).NumberFormat = $origFmt
    }
}

# --- Update data values (Receipts, Cost, Number of Plants, Year-to-Date figures) ---
$ws.Range("B8").Value = 1082182
$ws.Range("C8").Value = 1170593
$ws.Range("D8").Value = 2.09
$ws.Range("F8").Value = 279
$ws.Range("G8").Value = 340
$ws.Range("H8").Value = 11242072
$ws.Range("I8").Value = 13961955
$ws.Range("K8").Value = 2.23

$ws.Range("B9").Value = 9430
$ws.Range("C9").Value = 14148
$ws.Range("D9").Value = 10.07
$ws.Range("E9").Value = 8.96
$ws.Range("F9").Value = 161
$ws.Range("G9").Value = 194
$ws.Range("H9").Value = 92662
$ws.Range("I9").Value = 137610
$ws.Range("J9").Value = 9.24
$ws.Range("K9").Value = 11.68

$ws.Range("B10").Value = 9364
$ws.Range("C10").Value = 12117
$ws.Range("D10").Value = 2.26
$ws.Range("E10").Value = 1.59
$ws.Range("F10").Value = 9
$ws.Range("G10").Value = 12
$ws.Range("H10").Value = 106785
$ws.Range("I10").Value = 127631
$ws.Range("J10").Value = 1.61
$ws.Range("K10").Value = 1.86

$ws.Range("B11").Value = 721028
$ws.Range("C11").Value = 783337
$ws.Range("D11").Value = 3.02
$ws.Range("E11").Value = 2.65
$ws.Range("F11").Value = 744
$ws.Range("G11").Value = 763
$ws.Range("H11").Value = 9856886
$ws.Range("I11").Value = 9354902
$ws.Range("J11").Value = 2.79
$ws.Range("K11").Value = 3.29

$ws.Range("B12").Value = 1822005
$ws.Range("C12").Value = 1980194
$ws.Range("D12").Value = 2.47
$ws.Range("E12").Value = 2.38
$ws.Range("F12").Value = 922
$ws.Range("G12").Value = 977
$ws.Range("H12").Value = 21298406
$ws.Range("I12").Value = 23582098
$ws.Range("K12").Value = 2.67

$ws.Range("B19").Value = 801020
$ws.Range("C19").Value = 862786
$ws.Range("E19").Value = 2.2
$ws.Range("G19").Value = 223
$ws.Range("H19").Value = 8424882
$ws.Range("I19").Value = 10263092

$ws.Range("B20").Value = 6595
$ws.Range("C20").Value = 8558
$ws.Range("D20").Value = 9.79
$ws.Range("E20").Value = 8.8
$ws.Range("F20").Value = 94
$ws.Range("G20").Value = 120
$ws.Range("H20").Value = 67444
$ws.Range("I20").Value = 82638
$ws.Range("J20").Value = 9.03
$ws.Range("K20").Value = 11.57

$ws.Range("B21").Value = 7871
$ws.Range("C21").Value = 10082
$ws.Range("D21").Value = 2.22
$ws.Range("E21").Value = 1.46
$ws.Range("F21").Value = 7
$ws.Range("H21").Value = 91688
$ws.Range("I21").Value = 107437
$ws.Range("J21").Value = 1.48
$ws.Range("K21").Value = 1.8

$ws.Range("B22").Value = 338187
$ws.Range("C22").Value = 365361
$ws.Range("D22").Value = 3.37
$ws.Range("E22").Value = 2.97
$ws.Range("G22").Value = 406
$ws.Range("H22").Value = 4702060
$ws.Range("I22").Value = 4331629
$ws.Range("J22").Value = 3.08
$ws.Range("K22").Value = 3.57

$ws.Range("B23").Value = 1153673
$ws.Range("C23").Value = 1246786
$ws.Range("D23").Value = 2.54
$ws.Range("E23").Value = 2.47
$ws.Range("F23").Value = 528
$ws.Range("G23").Value = 543
$ws.Range("H23").Value = 13286075
$ws.Range("I23").Value = 14784797
$ws.Range("K23").Value = 2.69

$ws.Range("B30").Value = 270974
$ws.Range("C30").Value = 286023
$ws.Range("D30").Value = 1.92
$ws.Range("E30").Value = 1.97
$ws.Range("F30").Value = 67
$ws.Range("G30").Value = 87
$ws.Range("H30").Value = 2676066
$ws.Range("I30").Value = 3453390
$ws.Range("K30").Value = 2.11

$ws.Range("B31").Value = 2625
$ws.Range("C31").Value = 5410
$ws.Range("D31").Value = 10.76
$ws.Range("E31").Value = 9.13
$ws.Range("F31").Value = 55
$ws.Range("G31").Value = 60
$ws.Range("H31").Value = 23506
$ws.Range("I31").Value = 52848
$ws.Range("J31").Value = 9.85
$ws.Range("K31").Value = 11.78

$ws.Range("B32").Value = 1294
$ws.Range("C32").Value = 1643
$ws.Range("H32").Value = 12066
$ws.Range("I32").Value = 12809
$ws.Range("K32").Value = 2.45

$ws.Range("B33").Value = 317208
$ws.Range("C33").Value = 351912
$ws.Range("D33").Value = 2.6
$ws.Range("E33").Value = 2.31
$ws.Range("F33").Value = 292
$ws.Range("G33").Value = 301
$ws.Range("H33").Value = 4462966
$ws.Range("I33").Value = 4320981
$ws.Range("J33").Value = 2.45
$ws.Range("K33").Value = 3

$ws.Range("B34").Value = 592100
$ws.Range("C34").Value = 644987
$ws.Range("F34").Value = 343
$ws.Range("G34").Value = 371
$ws.Range("H34").Value = 7174604
$ws.Range("I34").Value = 7840027

$ws.Range("B41").Value = 237
$ws.Range("C41").Value = 182
$ws.Range("H41").Value = 1074
$ws.Range("I41").Value = 2252
$ws.Range("K41").Value = 2.86

$ws.Range("B44").Value = 613
$ws.Range("C44").Value = 775
$ws.Range("H44").Value = 7438
$ws.Range("I44").Value = 5975

$ws.Range("B45").Value = 850
$ws.Range("C45").Value = 957
$ws.Range("H45").Value = 8512
$ws.Range("I45").Value = 8226

$ws.Range("B52").Value = 9951
$ws.Range("C52").Value = 21602
$ws.Range("F52").Value = 17
$ws.Range("G52").Value = 29
$ws.Range("H52").Value = 140050
$ws.Range("I52").Value = 243222
$ws.Range("K52").Value = 2.72

$ws.Range("B53").Value = 210
$ws.Range("C53").Value = 180
$ws.Range("D53").Value = 10.4
$ws.Range("E53").Value = 11.49
$ws.Range("F53").Value = 12
$ws.Range("G53").Value = 14
$ws.Range("H53").Value = 1712
$ws.Range("I53").Value = 2124
$ws.Range("J53").Value = 9.74
$ws.Range("K53").Value = 13.64

$ws.Range("B54").Value = 200
$ws.Range("C54").Value = 393
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 3031
$ws.Range("I54").Value = 7385

$ws.Range("B55").Value = 65021
$ws.Range("C55").Value = 65289
$ws.Range("F55").Value = 46
$ws.Range("H55").Value = 684423
$ws.Range("I55").Value = 696317

$ws.Range("B56").Value = 75381
$ws.Range("C56").Value = 87464
$ws.Range("H56").Value = 829215
$ws.Range("I56").Value = 949048
